$d = $word.ActiveDocument

$pairs = @(
    @{old="2024-03-07 Thursday"; new="2024-03-08 Friday"},
    @{old="64÷8=8, 0"; new="56÷9=6, 2"},
    @{old="66÷2=33, 0"; new="81÷8=10, 1"},
    @{old="52÷2=26, 0"; new="95÷5=19, 0"},
    @{old="97÷6=16, 1"; new="33÷6=5, 3"},
    @{old="24÷4=6, 0"; new="78÷8=9, 6"},
    @{old="38÷3=12, 2"; new="80÷8=10, 0"},
    @{old="21÷4=5, 1"; new="45÷5=9, 0"},
    @{old="75÷7=10, 5"; new="60÷2=30, 0"},
    @{old="50÷3=16, 2"; new="90÷3=30, 0"},
    @{old="34÷3=11, 1"; new="91÷7=13, 0"},
    @{old="74÷4=18, 2"; new="39÷8=4, 7"},
    @{old="57÷9=6, 3"; new="59÷5=11, 4"},
    @{old="58÷7=8, 2"; new="81÷3=27, 0"},
    @{old="87÷7=12, 3"; new="88÷7=12, 4"},
    @{old="61÷5=12, 1"; new="55÷2=27, 1"},
    @{old="82÷9=9, 1"; new="52÷7=7, 3"},
    @{old="78÷5=15, 3"; new="96÷8=12, 0"},
    @{old="76÷3=25, 1"; new="63÷9=7, 0"},
    @{old="49÷9=5, 4"; new="47÷5=9, 2"},
    @{old="81÷9=9, 0"; new="84÷5=16, 4"},
    @{old="67÷2=33, 1"; new="36÷7=5, 1"},
    @{old="45÷6=7, 3"; new="68÷3=22, 2"},
    @{old="93÷5=18, 3"; new="63÷8=7, 7"},
    @{old="90÷9=10, 0"; new="23÷2=11, 1"},
    @{old="24÷7=3, 3"; new="29÷7=4, 1"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
